# Actualización automática 2025-09-11 11:20:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": row 38 (PORCEKER S.A.) gets new sales figures ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("C38").Value = 518.4
$wsGrupo.Range("E38").Value = 362.89
$wsGrupo.Range("M38").Value = 3668.44

# Totals row (57) counts of non-zero clients per column bump from 2 to 3
$wsGrupo.Range("C57").Value = "3 de 55"
$wsGrupo.Range("E57").Value = "3 de 55"

# --- Sheet "VENTA MENSUAL": septiembre column (F) for the same client/row ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F38").Value = 5945.84
$wsMensual.Range("F57").Value = 16238.94

# --- Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO per group ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2: 240X120 PORCELANATO
$wsCumpl.Range("D2").Value = 1933.63
$wsCumpl.Range("E2").Value = 958.57588040374
$wsCumpl.Range("F2").Value = 0.6685658213688692

# Row 4: FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 484.33
$wsCumpl.Range("E4").Value = 298.087163948959
$wsCumpl.Range("F4").Value = 0.6190176063565948

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 13789.1
$wsCumpl.Range("E12").Value = 33345.1631579098
$wsCumpl.Range("F12").Value = 0.2925493913801853

# Row 15: TOTAL
$wsCumpl.Range("D15").Value = 21080.74
$wsCumpl.Range("E15").Value = 77817.25992509275
$wsCumpl.Range("F15").Value = 0.2131563835059047

# Column width auto-adjustment on CUMPLIMIENTO MENSUAL (D and E got one unit narrower).
# ColumnWidth round-trips through this host with a fixed +5/6 character offset when
# written back to the OOXML <col width> attribute, so pre-compensate to land on the
# exact target stored width (13 and 22).
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - (5/6)
$wsCumpl.Columns.Item(5).ColumnWidth = 22 - (5/6)
